$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark newly-finished character assets as "final" (instead of "placeholder"),
# with their Source set to "original":
#  - ch-spouse-f  (row 3)
#  - ch-neighbor  (row 6)  -> old neighbour (peggy)
#  - ch-officer-2 (row 11) -> officer 2
$ws.Range("D3").Value = "final"
$ws.Range("E3").Value = "original"

$ws.Range("D6").Value = "final"
$ws.Range("E6").Value = "original"

$ws.Range("D11").Value = "final"
$ws.Range("E11").Value = "original"

# Leave the selection where the author ended up working
$ws.Range("F9").Select()
